$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1  = 0.41603216931191866
    2  = -0.0099999997604491853
    3  = -0.0089999997539074172
    4  = -0.011999999931362026
    5  = -0.0059999997579804898
    6  = -0.005999999751594487
    7  = -0.019999999705488491
    8  = -0.019999999705315297
    9  = -0.0059999997520669979
    10 = -0.0059999997533424221
    11 = -0.0044999997583019535
    12 = 0.076103156626948998
    13 = -0.005999999750101459
    14 = -0.011999999729194855
    15 = -0.0059999997479742717
    16 = -0.0059999997471233968
    17 = -0.0059999997459856402
    18 = -0.0089999997358711781
    19 = -0.02939299793582606
    20 = -0.0089999997613716687
    21 = -0.00899999976101995
    22 = -0.0089999997607748128
    23 = -0.058602940278441373
    24 = -0.082120671537858847
    25 = -0.041999999632198914
    26 = -0.0059999997508590752
    27 = -0.0059999997493922486
    28 = -0.0059999997436639418
    29 = -0.011999999720096355
    30 = -0.019999999691934889
    31 = -0.014999999705553435
    32 = -0.020999999685407111
    33 = -0.0059999997344588607
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
